$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders the weekly price-report rows (3-16, excluding the
# unchanged row 14) so that each row ends up with the full A:R content that
# used to belong to a different row - i.e. a pure row permutation, no new
# values. Capture each source row's full A:R content first (so overwrites
# don't clobber a value we still need to read later), then write them back
# out in the new order.

$rows = 3,4,5,6,7,8,9,10,11,12,13,14,15,16

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = $ws.Range("A$r`:R$r").Value()
}

# Maps each destination row to the source row whose data it should now hold.
$mapping = @{
    3  = 6
    4  = 7
    5  = 12
    6  = 16
    7  = 10
    8  = 11
    9  = 15
    10 = 9
    11 = 13
    12 = 4
    13 = 5
    14 = 14
    15 = 8
    16 = 3
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $ws.Range("A$destRow`:R$destRow").Value = $orig[$srcRow]
}
